$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 5458
$ws.Range("J2").Value = 6996.3335
$ws.Range("L2").Value = 6996.3335
$ws.Range("N2").Value = -7222.3335
$ws.Range("H9").Value = 3487.4614
$ws.Range("I9").Value = 4043.9
$ws.Range("J9").Value = 1632.6666
$ws.Range("K9").Value = 4043.9
$ws.Range("L9").Value = 1632.6666
$ws.Range("M9").Value = -3874.9
$ws.Range("N9").Value = -1970.6666
$ws.Range("H19").Value = 1436.8462
$ws.Range("J19").Value = 1216.5
$ws.Range("L19").Value = 1216.5
$ws.Range("N19").Value = -1566.5
$ws.Range("H33").Value = 338.125
$ws.Range("I33").Value = 295.83334
$ws.Range("K33").Value = 295.83334
$ws.Range("M33").Value = -66.83334000000002
$ws.Range("H112").Value = 1031.05
$ws.Range("J112").Value = 1006.3684
$ws.Range("L112").Value = 3019.1052
$ws.Range("N112").Value = -5235.1052
$ws.Range("H129").Value = 3342.25
$ws.Range("J129").Value = 3478.5557
$ws.Range("L129").Value = 10435.6671
$ws.Range("N129").Value = -20435.6671
$ws.Range("H138").Value = 3292.8823
$ws.Range("I138").Value = 3022.1667
$ws.Range("J138").Value = 3440.5454
$ws.Range("K138").Value = 9066.500100000001
$ws.Range("L138").Value = 10321.6362
$ws.Range("M138").Value = -3926.500100000001
$ws.Range("N138").Value = -20601.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1896.8334
$ws.Range("I45").Value = 1896.8334
$ws.Range("K45").Value = 1896.8334
$ws.Range("M45").Value = -1519.8334
$ws.Range("H55").Value = 176499.5
$ws.Range("J55").Value = 176499.5
$ws.Range("L55").Value = 176499.5
$ws.Range("N55").Value = -177129.5
$ws.Range("H80").Value = 103332.664
$ws.Range("J80").Value = 103332.664
$ws.Range("L80").Value = 103332.664
$ws.Range("N80").Value = -105328.664
$ws.Range("H83").Value = 103332.664
$ws.Range("J83").Value = 103332.664
$ws.Range("L83").Value = 309997.992
$ws.Range("N83").Value = -319981.992
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 22974
$ws.Range("I82").Value = 5632
$ws.Range("K82").Value = 5632
$ws.Range("M82").Value = -5249
$ws.Range("H85").Value = 22974
$ws.Range("I85").Value = 5632
$ws.Range("K85").Value = 5632
$ws.Range("M85").Value = -4306
$ws.Range("H105").Value = 3183.889
$ws.Range("I105").Value = 2957.5
$ws.Range("J105").Value = 4995
$ws.Range("K105").Value = 2957.5
$ws.Range("L105").Value = 4995
$ws.Range("M105").Value = -1210.5
$ws.Range("N105").Value = -8489
$ws.Range("H134").Value = 561.3333
$ws.Range("I134").Value = 561.3333
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 1683.9999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 851.0001
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H129").Value = 19000
$ws.Range("I129").Value = 19000
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 19000
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -14000
$ws.Range("H130").Value = 69989.336
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 69989.336
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 69989.336
$ws.Range("N130").Value = -80029.336
$ws.Range("H131").Value = 99999
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 99999
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 99999
$ws.Range("N131").Value = -110079
$ws.Range("H132").Value = 1385.95
$ws.Range("I132").Value = 936.1429000000001
$ws.Range("J132").Value = 2435.5
$ws.Range("K132").Value = 2808.4287
$ws.Range("L132").Value = 7306.5
$ws.Range("M132").Value = -278.4287000000004
$ws.Range("N132").Value = -12366.5
$ws.Range("H133").Value = 59162.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 59162.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 59162.5
$ws.Range("N133").Value = -64222.5
$ws.Range("H134").Value = 2786.3333
$ws.Range("I134").Value = 2668.6428
$ws.Range("J134").Value = 3198.25
$ws.Range("K134").Value = 8005.928400000001
$ws.Range("L134").Value = 9594.75
$ws.Range("M134").Value = -5470.928400000001
$ws.Range("N134").Value = -14664.75
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 466383.25
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 466383.25
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 466383.25
$ws.Range("N141").Value = -476743.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 12054
$ws.Range("I70").Value = 2466.5
$ws.Range("J70").Value = 15249.833
$ws.Range("K70").Value = 7399.5
$ws.Range("L70").Value = 45749.499
$ws.Range("M70").Value = -7084.5
$ws.Range("N70").Value = -46379.499
$ws.Range("H73").Value = 12054
$ws.Range("I73").Value = 2466.5
$ws.Range("J73").Value = 15249.833
$ws.Range("K73").Value = 7399.5
$ws.Range("L73").Value = 45749.499
$ws.Range("M73").Value = -6307.5
$ws.Range("N73").Value = -47933.499
$ws.Range("H134").Value = 15036.667
$ws.Range("J134").Value = 18992
$ws.Range("L134").Value = 56976
$ws.Range("N134").Value = -67116

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 53547
$ws.Range("J95").Value = 53547
$ws.Range("L95").Value = 53547
$ws.Range("N95").Value = -59039
$ws.Range("H122").Value = 1351.091
$ws.Range("I122").Value = 1351.091
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4053.273
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1603.273
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 40774.332
$ws.Range("J136").Value = 40774.332
$ws.Range("L136").Value = 122322.996
$ws.Range("N136").Value = -127422.996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H96").Value = 7108.4287
$ws.Range("J96").Value = 5495
$ws.Range("L96").Value = 5495
$ws.Range("N96").Value = -8241
